$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 86 - new TA/CRM sample reading, continuing the pattern of rows 77:85
$ws.Range("A85").Copy($ws.Range("A86"))
$ws.Range("A86").Value = 43820
$ws.Range("B86").Value = 2221.4069258804502
$ws.Range("C86").Value = 2207.0300000000002
$ws.Range("D86").Formula = "=100*(B86-C86)/C86"
$ws.Range("E86").Value = 169
$ws.Range("F86").Value = "New CRM opened 12/21/2019"

# Row 87 - CRM value only, start of next reading
$ws.Range("C87").Value = 2207.0300000000002

# Update view: scroll position + selection to match where the new data was entered
$ws.Range("C87").Select()
$excel.ActiveWindow.ScrollRow = 70
